# Adds the "ChinaMutualFundStockPortfolio" (中国共同基金股票持仓明细) table
# description rows to the FactorInfo sheet, and the matching table-registry
# row to the TableInfo sheet ("add fund stock pos").

$wb = $excel.ActiveWorkbook
$wsFactor = $wb.Worksheets.Item("FactorInfo")
$wsTable  = $wb.Worksheets.Item("TableInfo")

$rows = @(
    @{Row=875; A='中国共同基金股票持仓明细'; B='更新时间'; C='opdate'; D=0; E='None'},
    @{Row=876; A='中国共同基金股票持仓明细'; B='基金Wind代码'; C='s_info_windcode'; D=1; E='IDs'},
    @{Row=877; A='中国共同基金股票持仓明细'; B='截止日期'; C='f_prt_enddate'; D=1; E='date'},
    @{Row=878; A='中国共同基金股票持仓明细'; B='公告日期'; C='ann_date'; D=1; E='ann_dt'},
    @{Row=879; A='中国共同基金股票持仓明细'; B='持有股票Wind代码'; C='s_info_stockwindcode'; D=1; E='stock_id'},
    @{Row=880; A='中国共同基金股票持仓明细'; B='持有股票市值(元)'; C='f_prt_stkvalue'; D=0; E='None'},
    @{Row=881; A='中国共同基金股票持仓明细'; B='持有股票数量(股)'; C='f_prt_stkquantity'; D=0; E='None'},
    @{Row=882; A='中国共同基金股票持仓明细'; B='持有股票市值占基金净值比例(%)'; C='f_prt_stkvaluetonav'; D=0; E='None'},
    @{Row=883; A='中国共同基金股票持仓明细'; B='积极投资持有股票市值(元)'; C='f_prt_posstkvalue'; D=0; E='None'},
    @{Row=884; A='中国共同基金股票持仓明细'; B='积极投资持有股数(股)'; C='f_prt_posstkquantity'; D=0; E='None'},
    @{Row=885; A='中国共同基金股票持仓明细'; B='积极投资持有股票市值占净资产比例(%)'; C='f_prt_posstktonav'; D=0; E='None'},
    @{Row=886; A='中国共同基金股票持仓明细'; B='指数投资持有股票市值(元)'; C='f_prt_passtkevalue'; D=0; E='None'},
    @{Row=887; A='中国共同基金股票持仓明细'; B='指数投资持有股数(股)'; C='f_prt_passtkquantity'; D=0; E='None'},
    @{Row=888; A='中国共同基金股票持仓明细'; B='指数投资持有股票市值占净资产比例(%)'; C='f_prt_passtktonav'; D=0; E='None'},
    @{Row=889; A='中国共同基金股票持仓明细'; B='占股票市值比'; C='stock_per'; D=0; E='None'},
    @{Row=890; A='中国共同基金股票持仓明细'; B='占流通股本比'; C='float_shr_per'; D=0; E='None'}
)

foreach ($r in $rows) {
    $wsFactor.Range("A$($r.Row)").Value = $r.A
    $wsFactor.Range("B$($r.Row)").Value = $r.B
    $wsFactor.Range("C$($r.Row)").Value = $r.C
    $wsFactor.Range("D$($r.Row)").Value = $r.D
    $wsFactor.Range("E$($r.Row)").Value = $r.E
}

# Column E in this sheet carries a wrap/vertical-center style (the same one
# already used on E875:E877's placeholder cells). Brand-new rows (878:890)
# don't inherit it automatically, so stamp the format from an existing
# formatted cell onto the whole newly written E875:E890 block.
$wsFactor.Range("E874").Copy()
$wsFactor.Range("E875:E890").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTable.Range("A26").Value = '中国共同基金股票持仓明细'
$wsTable.Range("B26").Value = 'ChinaMutualFundStockPortfolio'
$wsTable.Range("C26").Value = '基金'

# Restore the selection the author left the workbook in after scrolling
# down to the newly-added rows (frozen-pane scroll position itself isn't
# part of the exposed object model, so only the active-cell selection is
# reproduced here).
[void]$wsFactor.Range("C887").Select()
[void]$wsTable.Range("B17").Select()
[void]$wsFactor.Select()
